$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "44.919.00"
Set-TextValue $ws.Range("E2") "  +0.51%  "
Set-TextValue $ws.Range("D3") "2.266.76"
Set-TextValue $ws.Range("E3") "  +1.30%  "
Set-TextValue $ws.Range("E4") "  -0.80%  "
Set-TextValue $ws.Range("D5") "302.10"
Set-TextValue $ws.Range("E5") "  -1.13%  "
Set-TextValue $ws.Range("D6") "94.29"
Set-TextValue $ws.Range("E6") "  -0.79%  "
Set-TextValue $ws.Range("D7") "0.565"
Set-TextValue $ws.Range("E7") "  -1.28%  "
Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.63%  "
Set-TextValue $ws.Range("D9") "0.511"
Set-TextValue $ws.Range("E9") "  -1.47%  "
Set-TextValue $ws.Range("D10") "34.31"
Set-TextValue $ws.Range("E10") "  -1.78%  "
Set-TextValue $ws.Range("D11") "0.0791"
Set-TextValue $ws.Range("E11") "  -1.44%  "
Set-TextValue $ws.Range("D12") "7.21"
Set-TextValue $ws.Range("E12") "  +0.39%  "
Set-TextValue $ws.Range("E13") "  -1.15%  "
Set-TextValue $ws.Range("D14") "2.606.49"
Set-TextValue $ws.Range("E14") "  +1.04%  "
Set-TextValue $ws.Range("D15") "2.267.54"
Set-TextValue $ws.Range("E15") "  +0.97%  "
Set-TextValue $ws.Range("D16") "13.61"
Set-TextValue $ws.Range("E16") "  +0.53%  "
Set-TextValue $ws.Range("D17") "0.801"
Set-TextValue $ws.Range("E17") "  -4.17%  "
Set-TextValue $ws.Range("D18") "44.817.26"
Set-TextValue $ws.Range("E18") "  +0.75%  "
Set-TextValue $ws.Range("D19") "13.05"
Set-TextValue $ws.Range("E19") "  +9.34%  "
Set-TextValue $ws.Range("D20") "0.0₃0923"
Set-TextValue $ws.Range("E20") "  -2.24%  "
Set-TextValue $ws.Range("D21") "6.05"
Set-TextValue $ws.Range("E21") "  -3.28%  "
Set-TextValue $ws.Range("D22") "65.66"
Set-TextValue $ws.Range("E22") "  +0.57%  "
Set-TextValue $ws.Range("D23") "238.50"
Set-TextValue $ws.Range("E23") "  -0.38%  "
Set-TextValue $ws.Range("E24") "  -2.14%  "
Set-TextValue $ws.Range("D25") "1.00"
Set-TextValue $ws.Range("E25") "  -0.29%  "
Set-TextValue $ws.Range("E26") "  -4.48%  "
Set-TextValue $ws.Range("D27") "41.36"
Set-TextValue $ws.Range("E27") "  +9.99%  "
Set-TextValue $ws.Range("D28") "2.30"
Set-TextValue $ws.Range("E28") "  -0.08%  "
Set-TextValue $ws.Range("E29") "  -1.95%  "
Set-TextValue $ws.Range("D30") "19.61"
Set-TextValue $ws.Range("D31") "152.61"
Set-TextValue $ws.Range("E31") "  +1.32%  "
Set-TextValue $ws.Range("D32") "5.56"
Set-TextValue $ws.Range("E32") "  -7.24%  "
Set-TextValue $ws.Range("D33") "0.0791"
Set-TextValue $ws.Range("E33") "  -0.40%  "
Set-TextValue $ws.Range("D34") "2.56"
Set-TextValue $ws.Range("E34") "  -2.65%  "
Set-TextValue $ws.Range("D35") "2.95"
Set-TextValue $ws.Range("E35") "  -2.78%  "
Set-TextValue $ws.Range("E36") "  -1.17%  "
Set-TextValue $ws.Range("E37") "  -3.18%  "
Set-TextValue $ws.Range("D38") "1.77"
Set-TextValue $ws.Range("E38") "  -3.77%  "
Set-TextValue $ws.Range("D39") "3.99"
Set-TextValue $ws.Range("E39") "  +6.41%  "
Set-TextValue $ws.Range("D40") "0.0310"
Set-TextValue $ws.Range("E40") "  +2.85%  "
Set-TextValue $ws.Range("D41") "3.23"
Set-TextValue $ws.Range("E41") "  -4.11%  "
Set-TextValue $ws.Range("D42") "13.65"
Set-TextValue $ws.Range("E42") "  -9.79%  "
Set-TextValue $ws.Range("D43") "0.998"
Set-TextValue $ws.Range("E43") "  -1.01%  "
Set-TextValue $ws.Range("D44") "1.91"
Set-TextValue $ws.Range("E44") "  +11.03%  "
Set-TextValue $ws.Range("D45") "1.739.96"
Set-TextValue $ws.Range("E45") "  -4.83%  "
Set-TextValue $ws.Range("D46") "0.194"
Set-TextValue $ws.Range("E46") "  +2.58%  "
Set-TextValue $ws.Range("D47") "76.47"
Set-TextValue $ws.Range("E47") "  -4.21%  "
Set-TextValue $ws.Range("D48") "69.48"
Set-TextValue $ws.Range("E48") "  +0.79%  "
Set-TextValue $ws.Range("E49") "  -3.24%  "
Set-TextValue $ws.Range("D50") "53.57"
Set-TextValue $ws.Range("E50") "  -1.24%  "
Set-TextValue $ws.Range("D51") "4.69"
Set-TextValue $ws.Range("E51") "  -3.74%  "
